$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$data = @(
    (2,8,116),
    (2,9,116),
    (2,10,0),
    (2,11,116),
    (2,12,0),
    (2,13,-3),
    (2,14,$null),
    (4,8,2378.8333),
    (4,9,759),
    (4,10,3998.6667),
    (4,11,759),
    (4,12,3998.6667),
    (4,13,-645),
    (4,14,-4226.6667),
    (9,8,597.73334),
    (9,9,269.72726),
    (9,10,1499.75),
    (9,11,269.72726),
    (9,12,1499.75),
    (9,13,-100.72726),
    (9,14,-1837.75),
    (19,8,914.95),
    (19,9,800),
    (19,10,1087.375),
    (19,11,800),
    (19,12,1087.375),
    (19,13,-625),
    (19,14,-1437.375),
    (28,8,61600.94),
    (28,10,94892.09),
    (28,12,94892.09),
    (28,14,-95862.09),
    (74,8,6352.0527),
    (74,9,6605.5625),
    (74,11,6605.5625),
    (74,13,-5669.5625),
    (77,8,6352.0527),
    (77,9,6605.5625),
    (77,11,33027.8125),
    (77,13,-28347.8125),
    (100,8,3083.6924),
    (100,9,3076.4443),
    (100,11,3076.4443),
    (100,13,-2535.4443),
    (111,8,52498.81),
    (111,10,99325.63),
    (111,12,297976.89),
    (111,14,-304110.89),
    (113,8,4662.4614),
    (113,9,4077.5715),
    (113,10,5344.8335),
    (113,11,4077.5715),
    (113,12,5344.8335),
    (113,13,-823.5715),
    (113,14,-11852.8335),
    (114,8,199500),
    (114,10,199500),
    (114,12,199500),
    (114,14,-208178),
    (115,8,500),
    (115,9,800),
    (115,10,200),
    (115,11,2400),
    (115,12,600),
    (115,13,-833),
    (115,14,-3734),
    (116,8,9567.75),
    (116,9,8453.333000000001),
    (116,11,8453.333000000001),
    (116,13,-5011.333000000001),
    (117,8,178583.33),
    (117,10,178583.33),
    (117,12,178583.33),
    (117,14,-187761.33),
    (118,8,2008.6666),
    (118,9,1982.125),
    (118,10,2221),
    (118,11,5946.375),
    (118,12,6663),
    (118,13,-4289.375),
    (118,14,-9977),
    (132,8,6173.5),
    (132,9,5898.1665),
    (132,10,6999.5),
    (132,11,17694.4995),
    (132,12,20998.5),
    (132,13,-15164.4995),
    (132,14,-26058.5),
    (135,8,2066.923),
    (135,9,460),
    (135,10,4638),
    (135,11,4140),
    (135,12,41742),
    (135,13,-1605),
    (135,14,-46812),
    (138,8,5463.8066),
    (138,9,4443.143),
    (138,10,7607.2),
    (138,11,13329.429),
    (138,12,22821.6),
    (138,13,-8189.429),
    (138,14,-33101.6),
    (141,8,8682.333000000001),
    (141,9,7547),
    (141,10,9250),
    (141,11,22641),
    (141,12,27750),
    (141,13,-17461),
    (141,14,-38110)
)
foreach ($t in $data) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$data = @(
    (2,8,567.1818),
    (2,9,567.1818),
    (2,11,567.1818),
    (2,13,-454.1818),
    (31,8,10041.667),
    (31,9,5675),
    (31,11,5675),
    (31,13,-5381),
    (32,8,5416.5713),
    (32,9,1814.5862),
    (32,11,1814.5862),
    (32,13,-1527.5862),
    (45,8,1646.5),
    (45,9,1004.4167),
    (45,11,1004.4167),
    (45,13,-627.4167),
    (61,8,21747038),
    (61,9,6168.6),
    (61,10,62511164),
    (61,11,6168.6),
    (61,12,62511164),
    (61,13,-5956.6),
    (61,14,-62511588),
    (63,8,2828.56),
    (63,9,1520.2727),
    (63,11,1520.2727),
    (63,13,-834.2727),
    (66,8,2828.56),
    (66,9,1520.2727),
    (66,11,7601.363499999999),
    (66,13,-4169.363499999999),
    (102,8,1318.1428),
    (102,9,1401.1111),
    (102,10,820.3333),
    (102,11,1401.1111),
    (102,12,820.3333),
    (102,13,220.8888999999999),
    (102,14,-4064.3333),
    (116,8,567.1818),
    (116,9,567.1818),
    (116,11,567.1818),
    (116,13,1726.8182),
    (117,8,199500),
    (117,10,199500),
    (117,12,199500),
    (117,14,-208678),
    (122,8,3361.2222),
    (122,9,3329.5293),
    (122,11,9988.5879),
    (122,13,-7538.5879),
    (136,8,21747038),
    (136,9,6168.6),
    (136,10,62511164),
    (136,11,18505.8),
    (136,12,187533492),
    (136,13,-15955.8),
    (136,14,-187538592)
)
foreach ($t in $data) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$data = @(
    (3,8,567.1818),
    (3,9,567.1818),
    (3,11,567.1818),
    (3,13,-453.1818),
    (94,8,666.931),
    (94,10,454.25),
    (94,12,454.25),
    (94,14,-1356.25),
    (99,8,1542.95),
    (99,9,1397.7778),
    (99,11,1397.7778),
    (99,13,100.2221999999999),
    (107,8,1509.8823),
    (107,9,880.2),
    (107,10,6232.5),
    (107,11,880.2),
    (107,12,6232.5),
    (107,13,1039.8),
    (107,14,-10072.5),
    (114,8,199500),
    (114,10,199500),
    (114,12,199500),
    (114,14,-208178),
    (115,8,199500),
    (115,10,199500),
    (115,12,199500),
    (115,14,-202634),
    (116,8,199500),
    (116,10,199500),
    (116,12,199500),
    (116,14,-208678),
    (117,8,199500),
    (117,10,199500),
    (117,12,199500),
    (117,14,-208678),
    (118,8,199500),
    (118,10,199500),
    (118,12,199500),
    (118,14,-202814),
    (139,8,32695),
    (139,9,23333.334),
    (139,10,60780),
    (139,11,23333.334),
    (139,12,60780),
    (139,13,-18193.334),
    (139,14,-71060)
)
foreach ($t in $data) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$data = @(
    (16,8,1242.8125),
    (16,9,1058.4482),
    (16,10,3025),
    (16,11,1058.4482),
    (16,12,3025),
    (16,13,-771.4482),
    (16,14,-3599),
    (22,8,633.41174),
    (22,9,465.45456),
    (22,10,941.3333),
    (22,11,465.45456),
    (22,12,941.3333),
    (22,13,-115.45456),
    (22,14,-1641.3333),
    (31,8,11469.375),
    (31,9,5399.8887),
    (31,10,19273),
    (31,11,5399.8887),
    (31,12,19273),
    (31,13,-5104.8887),
    (31,14,-19863),
    (34,8,11469.375),
    (34,9,5399.8887),
    (34,10,19273),
    (34,11,5399.8887),
    (34,12,19273),
    (34,13,-5197.8887),
    (34,14,-19677),
    (87,8,199500),
    (87,10,199500),
    (87,12,199500),
    (87,14,-201872),
    (90,8,199500),
    (90,10,199500),
    (90,12,598500),
    (90,14,-610356),
    (99,8,2448.508),
    (99,10,3099.1538),
    (99,12,3099.1538),
    (99,14,-6095.1538),
    (108,8,199500),
    (108,10,199500),
    (108,12,199500),
    (108,14,-207180),
    (110,8,197916.67),
    (110,10,197916.67),
    (110,12,197916.67),
    (110,14,-206096.67),
    (112,8,199357.14),
    (112,10,199357.14),
    (112,12,199357.14),
    (112,14,-202311.14),
    (113,8,1242.8125),
    (113,9,1058.4482),
    (113,10,3025),
    (113,11,1058.4482),
    (113,12,3025),
    (113,13,1111.5518),
    (113,14,-7365),
    (114,8,199500),
    (114,10,199500),
    (114,12,199500),
    (114,14,-208178),
    (115,8,159125),
    (115,10,199500),
    (115,12,199500),
    (115,14,-201850),
    (116,8,299500),
    (116,10,299500),
    (116,12,299500),
    (116,14,-308678),
    (117,8,169600),
    (117,9,50000),
    (117,10,199500),
    (117,11,50000),
    (117,12,199500),
    (117,13,-45411),
    (117,14,-208678),
    (122,8,1713.238),
    (122,9,1672.5264),
    (122,11,5017.5792),
    (122,13,-2567.5792),
    (126,8,2448.508),
    (126,10,3099.1538),
    (126,12,9297.4614),
    (126,14,-14237.4614),
    (132,8,5263.8423),
    (132,9,2866.6667),
    (132,11,8600.000100000001),
    (132,13,-6070.000100000001)
)
foreach ($t in $data) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$data = @(
    (2,8,763.7059),
    (2,9,122.9),
    (2,10,1679.1428),
    (2,11,737.4000000000001),
    (2,12,10074.8568),
    (2,13,-624.4000000000001),
    (2,14,-10300.8568),
    (10,8,1003.63635),
    (10,9,1103.8),
    (10,10,2),
    (10,11,3311.4),
    (10,12,6),
    (10,13,-3172.4),
    (10,14,-284),
    (12,8,286.86667),
    (12,10,392.6),
    (12,12,1177.8),
    (12,14,-1523.8),
    (17,8,1629.5714),
    (17,10,1000),
    (17,12,3000),
    (17,14,-3338),
    (21,8,330.66666),
    (21,9,296),
    (21,10,400),
    (21,11,888),
    (21,12,1200),
    (21,13,-715),
    (21,14,-1546),
    (23,8,135.33333),
    (23,9,120),
    (23,10,138.4),
    (23,11,360),
    (23,12,415.2),
    (23,13,-125),
    (23,14,-885.2),
    (25,8,3),
    (25,9,3),
    (25,11,9),
    (25,13,160),
    (29,8,252.25),
    (29,9,236.66667),
    (29,11,710.00001),
    (29,13,-433.00001),
    (30,8,3),
    (30,9,3),
    (30,11,9),
    (30,13,93),
    (31,8,2166),
    (31,9,500),
    (31,10,2999),
    (31,11,1500),
    (31,12,8997),
    (31,13,-1212),
    (31,14,-9573),
    (57,8,11499.777),
    (57,9,3499.6667),
    (57,11,10499.0001),
    (57,13,-9940.000100000001),
    (81,8,2142.8572),
    (81,10,1500),
    (81,12,4500),
    (81,14,-6746),
    (84,8,2142.8572),
    (84,10,1500),
    (84,12,13500),
    (84,14,-24732),
    (92,8,228),
    (92,9,137.5),
    (92,10,253.85715),
    (92,11,412.5),
    (92,12,761.5714499999999),
    (92,13,835.5),
    (92,14,-3257.57145),
    (97,8,1579.8),
    (97,9,3900),
    (97,10,999.75),
    (97,11,11700),
    (97,12,2999.25),
    (97,13,-11204),
    (97,14,-3991.25),
    (109,8,0),
    (109,9,0),
    (109,11,0),
    (109,13,$null),
    (113,8,1250.0588),
    (113,10,1250.0588),
    (113,12,3750.1764),
    (113,14,-8090.1764),
    (114,8,2245.6),
    (114,9,1676),
    (114,10,3100),
    (114,11,5028),
    (114,12,9300),
    (114,13,-1774),
    (114,14,-15808),
    (115,8,4990),
    (115,10,0),
    (115,12,0),
    (115,14,$null),
    (116,8,5000),
    (116,9,5000),
    (116,11,15000),
    (116,13,-11558),
    (117,8,892.625),
    (117,9,868.2857),
    (117,10,911.55554),
    (117,11,2604.8571),
    (117,12,2734.66662),
    (117,13,837.1428999999998),
    (117,14,-9618.66662),
    (118,8,5872.6665),
    (118,9,5872.6665),
    (118,11,17617.9995),
    (118,13,-16374.9995),
    (121,8,1407040.6),
    (121,9,619.8570999999999),
    (121,11,1859.5713),
    (121,13,-549.5712999999998),
    (122,8,600),
    (122,10,0),
    (122,12,0),
    (122,14,$null)
)
foreach ($t in $data) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$data = @(
    (5,8,42499.5),
    (5,9,42499.5),
    (5,10,0),
    (5,11,42499.5),
    (5,12,0),
    (5,13,-42387.5),
    (5,14,$null),
    (31,8,4965),
    (31,9,4958),
    (31,11,4958),
    (31,13,-4666),
    (37,8,4965),
    (37,9,4958),
    (37,11,4958),
    (37,13,-4681),
    (102,8,2916.6875),
    (102,9,2287.7856),
    (102,10,7319),
    (102,11,2287.7856),
    (102,12,7319),
    (102,13,-665.7856000000002),
    (102,14,-10563),
    (116,8,199500),
    (116,10,199500),
    (116,12,199500),
    (116,14,-208678),
    (117,8,197968.33),
    (117,10,197968.33),
    (117,12,197968.33),
    (117,14,-204852.33),
    (118,8,199500),
    (118,10,199500),
    (118,12,199500),
    (118,14,-202814),
    (122,8,7015.268),
    (122,9,5926.4165),
    (122,10,8552.471),
    (122,11,17779.2495),
    (122,12,25657.413),
    (122,13,-15329.2495),
    (122,14,-30557.413),
    (123,8,63950),
    (123,10,63950),
    (123,12,63950),
    (123,14,-68850),
    (132,8,6247.9165),
    (132,9,2996.75),
    (132,10,7873.5),
    (132,11,8990.25),
    (132,12,23620.5),
    (132,13,-6460.25),
    (132,14,-28680.5),
    (136,8,15850.3),
    (136,10,16833.666),
    (136,12,50500.99800000001),
    (136,14,-55600.99800000001)
)
foreach ($t in $data) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$data = @(
    (40,8,2885.9688),
    (40,9,2235.8572),
    (40,11,2235.8572),
    (40,13,-2099.8572),
    (46,8,3182.1667),
    (46,9,512),
    (46,10,4209.154),
    (46,11,512),
    (46,12,4209.154),
    (46,13,-324),
    (46,14,-4585.154),
    (55,8,1248.1538),
    (55,9,247.33333),
    (55,11,247.33333),
    (55,13,-74.33332999999999),
    (61,8,4409.65),
    (61,9,1926),
    (61,10,6893.3),
    (61,11,1926),
    (61,12,6893.3),
    (61,13,-1724),
    (61,14,-7297.3),
    (82,8,1601.4546),
    (82,9,1015),
    (82,11,1015),
    (82,13,-654),
    (85,8,1601.4546),
    (85,9,1015),
    (85,11,1015),
    (85,13,233),
    (93,8,3209),
    (93,9,3013.1428),
    (93,10,3666),
    (93,11,3013.1428),
    (93,12,3666),
    (93,13,-1765.1428),
    (93,14,-6162),
    (100,8,2086.889),
    (100,9,2086.889),
    (100,11,2086.889),
    (100,13,-1545.889),
    (109,8,156485),
    (109,10,156485),
    (109,12,156485),
    (109,14,-159259),
    (113,8,4409.65),
    (113,9,1926),
    (113,10,6893.3),
    (113,11,1926),
    (113,12,6893.3),
    (113,13,244),
    (113,14,-11233.3),
    (114,8,199500),
    (114,10,199500),
    (114,12,199500),
    (114,14,-208178),
    (117,8,199500),
    (117,10,199500),
    (117,12,199500),
    (117,14,-208678),
    (118,8,199500),
    (118,10,199500),
    (118,12,199500),
    (118,14,-202814),
    (122,8,5875.6924),
    (122,9,4884.857),
    (122,11,14654.571),
    (122,13,-12204.571),
    (132,8,9540.25),
    (132,9,7690.2856),
    (132,11,23070.8568),
    (132,13,-20540.8568),
    (136,8,38468428),
    (136,9,6611.875),
    (136,10,100007330),
    (136,11,19835.625),
    (136,12,300021990),
    (136,13,-17285.625),
    (136,14,-300027090)
)
foreach ($t in $data) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$data = @(
    (96,8,15774.2),
    (96,9,912.5),
    (96,10,21178.455),
    (96,11,912.5),
    (96,12,21178.455),
    (96,13,460.5),
    (96,14,-23924.455),
    (107,8,1491.3334),
    (107,9,816.7273),
    (107,10,2551.4285),
    (107,11,2450.1819),
    (107,12,7654.2855),
    (107,13,-530.1819),
    (107,14,-11494.2855),
    (112,8,136610.88),
    (112,10,136610.88),
    (112,12,136610.88),
    (112,14,-139564.88),
    (117,8,199500),
    (117,10,199500),
    (117,12,199500),
    (117,14,-208678),
    (118,8,199500),
    (118,10,199500),
    (118,12,199500),
    (118,14,-202814),
    (122,8,2627.6),
    (122,9,2566.353),
    (122,11,7699.059),
    (122,13,-5249.059),
    (126,8,2469.9333),
    (126,9,2220.75),
    (126,11,6662.25),
    (126,13,-4192.25),
    (132,8,2367.6428),
    (132,9,1804.9524),
    (132,10,4055.7144),
    (132,11,5414.857199999999),
    (132,12,12167.1432),
    (132,13,-2884.857199999999),
    (132,14,-17227.1432),
    (136,8,7745.7),
    (136,9,8042.8423),
    (136,11,24128.5269),
    (136,13,-21578.5269)
)
foreach ($t in $data) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}
